# Weekly update: insert a new record as the second entry for
# "Vega Modelo de Temuco - Espárragos" (row 85), pushing the existing
# rows 85-93 down to 86-94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 85 (shifts rows 85:93 down to 86:94,
# just like Excel's own Rows.Insert / xlShiftDown behaviour).
$ws.Rows("85:85").Insert()

# Populate the new row 85 with the latest weekly price record.
$ws.Cells.Item(85, 1).Value  = 10
$ws.Cells.Item(85, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(85, 3).Value  = "La Araucanía"
$ws.Cells.Item(85, 4).Value  = 45209
$ws.Cells.Item(85, 5).Value  = 9
$ws.Cells.Item(85, 6).Value  = 300000000
$ws.Cells.Item(85, 7).Value  = "Espárragos"
$ws.Cells.Item(85, 8).Value  = "Sin especificar"
$ws.Cells.Item(85, 9).Value  = "Primera"
$ws.Cells.Item(85, 10).Value = 300
$ws.Cells.Item(85, 11).Value = 1500
$ws.Cells.Item(85, 12).Value = 1500
$ws.Cells.Item(85, 13).Value = 1500
$ws.Cells.Item(85, 14).Value = "$/kilo"
$ws.Cells.Item(85, 15).Value = "Región del Maule"
$ws.Cells.Item(85, 16).Value = 1500
$ws.Cells.Item(85, 17).Value = 1
$ws.Cells.Item(85, 18).Value = "Hortaliza"
